{"js": "// Office.js (Word JavaScript API) edit script.\n// This is the body of `async (context) => { ... }`.\n//\n// Summary of the change (see commit message / diff):\n//   1. Insert a brand-new paragraph \"Opportunity number: XXXX-XXXX-XXXX\"\n//      (two runs) right before the existing \"Opportunity name: ...\" paragraph.\n//   2. In the \"Opportunity name: ...\" paragraph, merge the trailing \", t\" +\n//      \"esting fixture file\" runs into a single \", testing fixture file\" run.\n\nconst body = context.document.body;\n\n// --- 1. Locate the \"Opportunity name:\" paragraph -------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet opportunityNamePara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(\"Opportunity name:\") === 0) {\n    opportunityNamePara = p;\n    break;\n  }\n}\nif (!opportunityNamePara) {\n  throw new Error(\"Could not find the 'Opportunity name:' paragraph\");\n}\n\n// --- 2. Insert the new \"Opportunity number: XXXX-XXXX-XXXX\" paragraph ----\n// Office.js normally collapses adjacent runs that share identical\n// formatting, so a plain insertParagraph()/insertText() pair would collapse\n// back into a single run. Inserting raw OOXML (Flat OPC package) keeps the\n// two runs distinct, matching the target markup exactly. A trailing empty\n// paragraph is required so Word treats the fragment as whole paragraph(s)\n// rather than splicing the runs into the existing \"Opportunity name:\" text;\n// that stray empty paragraph is removed again afterwards.\nconst newParaOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  \"<w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\">Opportunity number: </w:t></w:r>' +\n  \"<w:r><w:t>XXXX-XXXX-XXXX</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"<w:p></w:p>\" +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nopportunityNamePara.insertOoxml(newParaOoxml, \"Start\");\nawait context.sync();\n\n// Remove the stray empty paragraph that insertOoxml left behind directly\n// above the (now shifted) \"Opportunity name:\" paragraph.\nconst paragraphsAfterInsert = body.paragraphs;\nparagraphsAfterInsert.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphsAfterInsert.items.length; i++) {\n  const p = paragraphsAfterInsert.items[i];\n  const next = paragraphsAfterInsert.items[i + 1];\n  if (p.text === \"\" && next && next.text.indexOf(\"Opportunity name:\") === 0) {\n    p.delete();\n    break;\n  }\n}\nawait context.sync();\n\n// --- 3. Merge the \", t\" + \"esting fixture file\" runs ---------------------\nconst startRun = body.search(\", t\", { matchCase: true }).getFirst();\nconst endRun = body.search(\"esting fixture file\", { matchCase: true }).getFirst();\nconst mergedRange = startRun.expandTo(endRun);\nmergedRange.insertText(\", testing fixture file\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Summary of the change (see commit message / diff):\n#   1. Insert a brand-new paragraph \"Opportunity number: XXXX-XXXX-XXXX\"\n#      (two runs) right before the existing \"Opportunity name: ...\" paragraph.\n#   2. In the \"Opportunity name: ...\" paragraph, merge the trailing \", t\" +\n#      \"esting fixture file\" runs into a single \", testing fixture file\" run.\n\n$d = $word.ActiveDocument\n\n# --- 1. Locate the \"Opportunity name:\" paragraph -------------------------\n$paras = $d.Paragraphs\n$targetIdx = -1\n$i = 0\nforeach ($p in $paras) {\n    $i = $i + 1\n    if ($p.Range.Text.StartsWith(\"Opportunity name:\")) {\n        $targetIdx = $i\n        break\n    }\n}\nif ($targetIdx -eq -1) {\n    throw \"Could not find the 'Opportunity name:' paragraph\"\n}\n\n$target = $d.Paragraphs.Item($targetIdx)\n$targetRange = $target.Range\n$targetRange.Collapse(1)          # wdCollapseStart\n$targetRange.InsertParagraphBefore()\n\n# --- 2. Fill the newly created (now empty) paragraph with two runs -------\n# A plain InsertAfter/InsertAfter pair gets canonicalized back into a single\n# run (same text formatting on both), so use InsertXML with a tiny\n# WordprocessingML fragment to keep \"Opportunity number: \" and\n# \"XXXX-XXXX-XXXX\" as two distinct runs, matching the target markup.\n$newPara = $d.Paragraphs.Item($targetIdx)\n$newParaRange = $newPara.Range\n$newParaXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:t xml:space=\"preserve\">Opportunity number: </w:t></w:r><w:r><w:t>XXXX-XXXX-XXXX</w:t></w:r></w:p>'\n$newParaRange.InsertXML($newParaXml)\n\n# --- 3. Merge the \", t\" + \"esting fixture file\" runs ---------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \", t\"\n$find.Execute() | Out-Null\n$startRange = $d.Content.Duplicate\n$startRange.Start = $find.Parent.Start\n$startRange.End = $find.Parent.End\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"esting fixture file\"\n$find2.Execute() | Out-Null\n$endRange = $d.Content.Duplicate\n$endRange.Start = $find2.Parent.Start\n$endRange.End = $find2.Parent.End\n\n$mergedRange = $d.Range($startRange.Start, $endRange.End)\n$mergedRange.Text = \", testing fixture file\"\n"}
